$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row 3 corresponds to b5fa3b00-7fc8-45f1-91f6-e18accf757cd.md
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-10-21 00:33:47"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to b5fa3b00-7fc8-45f1-91f6-e18accf757cd.md
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-10-21 00:33:35"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d72db1c23a1c9a38c6aa5d9505302e99f273040/e2e/b5fa3b00-7fc8-45f1-91f6-e18accf757cd.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/24f80cbc3320e2778c4c2f9c1f97fd893a46bf53/e2e/b5fa3b00-7fc8-45f1-91f6-e18accf757cd.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet: row 3 corresponds to b5fa3b00-7fc8-45f1-91f6-e18accf757cd.md
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-10-21 00:33:47"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d72db1c23a1c9a38c6aa5d9505302e99f273040/e2e/b5fa3b00-7fc8-45f1-91f6-e18accf757cd.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/24f80cbc3320e2778c4c2f9c1f97fd893a46bf53/e2e/b5fa3b00-7fc8-45f1-91f6-e18accf757cd.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
